$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C93").Value = "Drawdown_Deviation_test"
$ws.Range("B93").Value = "Test drawdown deviation"
$ws.Range("A93").Value = "Drawdown Deviation"

$ws.Range("B100").Select()
